$wb = $excel.ActiveWorkbook

# --- Sheet: LH_TC_NOTIFICATION_REVIEWS ---
$ws1 = $wb.Worksheets.Item("LH_TC_NOTIFICATION_REVIEWS")
$ws1.Range("I5").Value = "Open"
$ws1.Range("J5").Value = "Open"

# --- Sheet: Version History ---
$ws2 = $wb.Worksheets.Item("Version History")
$ws2.Rows.Item(6).Delete()

# Update selections / scroll positions to match final state
$ws1.Range("J5").Select()
$ws1.Application.ActiveWindow.ScrollColumn = 6

$ws2.Activate()
$ws2.Range("A16").Select()
